# Update countries & provincias Spain
# Applies updated COVID-19 case data for several countries, re-positions the
# Serbia and Tanzania rows (moving them right after Peru and Etiopia
# respectively, with freshly updated figures) while the rows that used to sit
# in between keep their original data and simply shift down by one row, and
# refreshes the "datos actualizados" timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value refreshes (no row movement) -----------------------------

# Reino Unido (row 11)
$ws.Cells.Item(11, 2).Value = 29474
$ws.Cells.Item(11, 3).Value = 4324
$ws.Cells.Item(11, 5).Value = 26987
$ws.Cells.Item(11, 7).Value = 563
$ws.Cells.Item(11, 8).Value = 2352

# Suiza (row 12)
$ws.Cells.Item(12, 2).Value = 17137
$ws.Cells.Item(12, 3).Value = 532
$ws.Cells.Item(12, 5).Value = 13709
$ws.Cells.Item(12, 6).Value = 348
$ws.Cells.Item(12, 7).Value = 28
$ws.Cells.Item(12, 8).Value = 461

# Bosnia y Herzegovina (row 72)
$ws.Cells.Item(72, 2).Value = 455
$ws.Cells.Item(72, 3).Value = 35
$ws.Cells.Item(72, 4).Value = 19
$ws.Cells.Item(72, 5).Value = 423

# Vietnam (row 93)
$ws.Cells.Item(93, 2).Value = 218
$ws.Cells.Item(93, 3).Value = 6
$ws.Cells.Item(93, 5).Value = 155

# --- Serbia moves from after Colombia to right after Peru -----------------
# Rows 50-53 (Argentina, Croacia, Singapur, Colombia) shift down one row to
# 51-54, keeping their existing data untouched. Row 50 becomes Serbia with
# refreshed totals.

$ws.Cells.Item(54, 1).Value = "Colombia"
$ws.Cells.Item(54, 2).Value = 906
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 31
$ws.Cells.Item(54, 5).Value = 859
$ws.Cells.Item(54, 6).Value = 35
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 16

$ws.Cells.Item(53, 1).Value = "Singapur"
$ws.Cells.Item(53, 2).Value = 926
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 240
$ws.Cells.Item(53, 5).Value = 683
$ws.Cells.Item(53, 6).Value = 22
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 3

$ws.Cells.Item(52, 1).Value = "Croacia"
$ws.Cells.Item(52, 2).Value = 963
$ws.Cells.Item(52, 3).Value = 96
$ws.Cells.Item(52, 4).Value = 73
$ws.Cells.Item(52, 5).Value = 884
$ws.Cells.Item(52, 6).Value = 34
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 6

$ws.Cells.Item(51, 1).Value = "Argentina"
$ws.Cells.Item(51, 2).Value = 1054
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 240
$ws.Cells.Item(51, 5).Value = 787
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 27

$ws.Cells.Item(50, 1).Value = "Serbia"
$ws.Cells.Item(50, 2).Value = 1060
$ws.Cells.Item(50, 3).Value = 160
$ws.Cells.Item(50, 4).Value = 42
$ws.Cells.Item(50, 5).Value = 990
$ws.Cells.Item(50, 6).Value = 62
$ws.Cells.Item(50, 7).Value = 5
$ws.Cells.Item(50, 8).Value = 28

# --- Tanzania moves from after Congo to right after Etiopia ---------------
# Row 149 (Congo) shifts down to row 150, keeping its existing data, and row
# 149 becomes Tanzania with refreshed totals.

$ws.Cells.Item(150, 1).Value = "Congo"
$ws.Cells.Item(150, 2).Value = 19
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 0
$ws.Cells.Item(150, 5).Value = 19
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 0

$ws.Cells.Item(149, 1).Value = "Tanzania"
$ws.Cells.Item(149, 2).Value = 20
$ws.Cells.Item(149, 3).Value = 1
$ws.Cells.Item(149, 4).Value = 1
$ws.Cells.Item(149, 5).Value = 18
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 1

# --- Refresh "last updated" banner -----------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 15:20"
